$wb = $excel.ActiveWorkbook

# The "transformer" sheet needs a new column inserted before column H,
# with header value "b", and it should become the active sheet with H1 selected.
$ws = $wb.Worksheets.Item("transformer")

# Insert a new column at H, shifting existing H:Q headers to I:R
$ws.Columns("H:H").Insert()

# Set the header text for the newly inserted column
$ws.Range("H1").Value = "b"

# Make this sheet the active one, with H1 selected
$ws.Activate()
$ws.Range("H1").Select()
